$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column in H1, matching the header style used by the
# other header cells (e.g. style used in G1).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats

# Add the Save column values for the data rows.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
